$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format temporarily so numeric-looking strings
# (e.g. "64.29") are stored as literal text instead of being auto-converted
# to numbers by Excel's normal type inference on assignment.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '28.834.00'
$ws.Range('E2').Value = '  +2.00%  '
$ws.Range('D3').Value = '1.576.71'
$ws.Range('E3').Value = '  +1.68%  '
$ws.Range('E4').Value = '  -0.62%  '
$ws.Range('D5').Value = '211.09'
$ws.Range('E5').Value = '  +0.93%  '
$ws.Range('D6').Value = '0.519'
$ws.Range('E6').Value = '  +7.07%  '
$ws.Range('D7').Value = '0.994'
$ws.Range('E7').Value = '  -0.78%  '
$ws.Range('D8').Value = '25.38'
$ws.Range('E8').Value = '  +8.68%  '
$ws.Range('D9').Value = '0.247'
$ws.Range('E9').Value = '  +2.33%  '
$ws.Range('D10').Value = '0.0591'
$ws.Range('E10').Value = '  +1.45%  '
$ws.Range('D11').Value = '0.0899'
$ws.Range('E11').Value = '  +1.21%  '
$ws.Range('D12').Value = '1.799.49'
$ws.Range('E12').Value = '  +1.55%  '
$ws.Range('D13').Value = '1.587.59'
$ws.Range('E13').Value = '  +2.41%  '
$ws.Range('D14').Value = '28.864.24'
$ws.Range('E14').Value = '  +2.13%  '
$ws.Range('D15').Value = '0.519'
$ws.Range('E15').Value = '  +2.35%  '
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('D17').Value = '62.10'
$ws.Range('E17').Value = '  +2.88%  '
$ws.Range('D18').Value = '231.68'
$ws.Range('E18').Value = '  +2.21%  '
$ws.Range('D19').Value = '7.36'
$ws.Range('E19').Value = '  +0.78%  '
$ws.Range('D20').Value = '0.0₃0690'
$ws.Range('E20').Value = '  +2.67%  '
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('D22').Value = '3.97'
$ws.Range('E22').Value = '  +1.73%  '
$ws.Range('D23').Value = '9.12'
$ws.Range('E23').Value = '  +3.81%  '
$ws.Range('E24').Value = '  +4.45%  '
$ws.Range('D25').Value = '152.13'
$ws.Range('E25').Value = '  +3.20%  '
$ws.Range('E26').Value = '  +4.41%  '
$ws.Range('D27').Value = '14.92'
$ws.Range('E27').Value = '  +1.26%  '
$ws.Range('E28').Value = '  +1.74%  '
$ws.Range('E29').Value = '  -0.60%  '
$ws.Range('D30').Value = '0.0462'
$ws.Range('E30').Value = '  -0.71%  '
$ws.Range('D31').Value = '1.06'
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('E32').Value = '  +1.15%  '
$ws.Range('D33').Value = '1.414.84'
$ws.Range('E33').Value = '  +2.37%  '
$ws.Range('E34').Value = '  -0.58%  '
$ws.Range('E35').Value = '  -1.70%  '
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('E37').Value = '  +7.43%  '
$ws.Range('E38').Value = '  -2.21%  '
$ws.Range('D39').Value = '0.0163'
$ws.Range('E39').Value = '  +1.15%  '
$ws.Range('D40').Value = '0.523'
$ws.Range('E40').Value = '  +2.96%  '
$ws.Range('E41').Value = '  +1.87%  '
$ws.Range('E42').Value = '  -0.58%  '
$ws.Range('D43').Value = '0.782'
$ws.Range('E43').Value = '  +1.53%  '
$ws.Range('E44').Value = '  -0.92%  '
$ws.Range('D45').Value = '64.29'
$ws.Range('E45').Value = '  +4.39%  '
$ws.Range('D46').Value = '5.28'
$ws.Range('E46').Value = '  -2.06%  '
$ws.Range('D47').Value = '1.712.77'
$ws.Range('E47').Value = '  +1.66%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '85.08'
$ws.Range('E48').Value = '  -0.25%  '
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').Value = '0.830'
$ws.Range('E49').Value = '  -8.68%  '
$ws.Range('D50').Value = '42.77'
$ws.Range('E50').Value = '  +2.91%  '
$ws.Range('E51').Value = '  +0.52%  '

# Restore column D to its original (unstyled / General) appearance so no
# stray style index is left attached to the cells.
$ws.Range("D2:D51").Style = "Normal"

